# "seperated game from vis"
# The "assets" sheet lists game asset file paths (column A) with helper
# formulas in columns B:F that derive shortened names / quoted JS literals.
# The epicon icon set was reworked (several icons renamed/retired, several
# new icons added) and the list re-sorted alphabetically, which grew the
# block from 27 to 29 rows. The "planet backgrounds" block further down the
# sheet was pushed down to make room (with a small gap of blank rows left
# behind, rows 80:86), moving from rows 78-86 to rows 87-95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: grow the epicon block by 2 rows (51-77 -> 51-79) and open
#    up a further 7-row gap before the "planet backgrounds" header so it
#    (and the spacebox rows that follow it) end up at rows 87-95.
# ---------------------------------------------------------------------
$ws.Rows("78:79").Insert()
$ws.Rows("80:86").Insert()

# ---------------------------------------------------------------------
# 2. Write the new, re-sorted epicon file list into column A (rows 51-79).
# ---------------------------------------------------------------------
$epiconPaths = @(
  "C:\py\mr-game-webapp\images\epicons\bmi.png",
  "C:\py\mr-game-webapp\images\epicons\caffeine.png",
  "C:\py\mr-game-webapp\images\epicons\chd_alt.png",
  "C:\py\mr-game-webapp\images\epicons\diabetes.png",
  "C:\py\mr-game-webapp\images\epicons\drinking.png",
  "C:\py\mr-game-webapp\images\epicons\drugs.png",
  "C:\py\mr-game-webapp\images\epicons\education_schoolYears.png",
  "C:\py\mr-game-webapp\images\epicons\education_schoolYears_alt.png",
  "C:\py\mr-game-webapp\images\epicons\eveningness.png",
  "C:\py\mr-game-webapp\images\epicons\exercise_bike.png",
  "C:\py\mr-game-webapp\images\epicons\gaming.png",
  "C:\py\mr-game-webapp\images\epicons\intelligence.png",
  "C:\py\mr-game-webapp\images\epicons\intelligence_alt.png",
  "C:\py\mr-game-webapp\images\epicons\mh_anxiety2.png",
  "C:\py\mr-game-webapp\images\epicons\mh_depression2.png",
  "C:\py\mr-game-webapp\images\epicons\mh_ocd2.png",
  "C:\py\mr-game-webapp\images\epicons\phone.png",
  "C:\py\mr-game-webapp\images\epicons\sleep.png",
  "C:\py\mr-game-webapp\images\epicons\sleep_duration.png",
  "C:\py\mr-game-webapp\images\epicons\sleep_insomnia.png",
  "C:\py\mr-game-webapp\images\epicons\smoking.png",
  "C:\py\mr-game-webapp\images\epicons\social_chatBubble.png",
  "C:\py\mr-game-webapp\images\epicons\social_chatBubble_alt.png",
  "C:\py\mr-game-webapp\images\epicons\social_loneliness.png",
  "C:\py\mr-game-webapp\images\epicons\social_loneliness_alt.png",
  "C:\py\mr-game-webapp\images\epicons\social_notChatting.png",
  "C:\py\mr-game-webapp\images\epicons\wellbeing2_alt.png",
  "C:\py\mr-game-webapp\images\epicons\work_nightShifts.png",
  "C:\py\mr-game-webapp\images\epicons\work_nightShifts_alt.png"
)

for ($i = 0; $i -lt $epiconPaths.Length; $i++) {
    $row = 51 + $i
    $ws.Range("A$row").Value = $epiconPaths[$i]
}

# ---------------------------------------------------------------------
# 3. Re-apply the derivation formulas for the whole, now-29-row, block.
#    (Columns B:F simply strip/re-quote the path in column A.)
# ---------------------------------------------------------------------
$ws.Range("B51:B79").FormulaR1C1 = "=RIGHT(RC[-1],LEN(RC[-1])-33)"
$ws.Range("C51:C79").FormulaR1C1 = "=LEFT(RIGHT(RC[-2],LEN(RC[-2])-33),LEN(RC[-2])-37)"
$ws.Range("D51:D79").FormulaR1C1 = "=""'""&RC[-1]&""',"""
$ws.Range("E51:E79").FormulaR1C1 = "=RIGHT(RC[-4],LEN(RC[-4])-21)"
$ws.Range("F51:F79").FormulaR1C1 = "=""'""&RC[-1]&""',"""

# ---------------------------------------------------------------------
# 4. Update the sheet selection to reflect where the editor was last
#    working (column F across the new epicon block).
# ---------------------------------------------------------------------
$ws.Range("F51:F79").Select()
$ws.Range("F51").Activate()
